$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 283.875
$ws.Range("I28").Value = 118.25
$ws.Range("K28").Value = 118.25
$ws.Range("M28").Value = 366.75

$ws.Range("H32").Value = 1012.8333
$ws.Range("I32").Value = 799
$ws.Range("J32").Value = 1119.75
$ws.Range("K32").Value = 799
$ws.Range("L32").Value = 1119.75
$ws.Range("M32").Value = -473
$ws.Range("N32").Value = -1771.75

$ws.Range("H33").Value = 277.14285
$ws.Range("I33").Value = 202.94118
$ws.Range("K33").Value = 202.94118
$ws.Range("M33").Value = 26.05882

$ws.Range("H40").Value = 6294.1177
$ws.Range("I40").Value = 3000
$ws.Range("K40").Value = 3000
$ws.Range("M40").Value = -2825

$ws.Range("H58").Value = 1514.2858
$ws.Range("J58").Value = 3333.3333
$ws.Range("L58").Value = 9999.999899999999
$ws.Range("N58").Value = -10299.9999

$ws.Range("H61").Value = 61.5
$ws.Range("I61").Value = 61.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 184.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -12.5
$ws.Range("N61").ClearContents()

$ws.Range("H76").Value = 1718.6666
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 1718.6666
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 1718.6666
$ws.Range("N76").Value = -2348.6666
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 1718.6666
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 1718.6666
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 1718.6666
$ws.Range("N79").Value = -3902.6666
$ws.Range("M79").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 5014.5
$ws.Range("I36").Value = 5014.5
$ws.Range("K36").Value = 5014.5
$ws.Range("M36").Value = -4668.5

$ws.Range("H63").Value = 991.5
$ws.Range("I63").Value = 1059.8
$ws.Range("J63").Value = 650
$ws.Range("K63").Value = 1059.8
$ws.Range("L63").Value = 650
$ws.Range("M63").Value = -373.8
$ws.Range("N63").Value = -2022

$ws.Range("H66").Value = 991.5
$ws.Range("I66").Value = 1059.8
$ws.Range("J66").Value = 650
$ws.Range("K66").Value = 5299
$ws.Range("L66").Value = 3250
$ws.Range("M66").Value = -1867
$ws.Range("N66").Value = -10114

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H102").Value = 1828.8889
$ws.Range("I102").Value = 1828.8889
$ws.Range("K102").Value = 1828.8889
$ws.Range("M102").Value = -206.8888999999999

$ws.Range("H110").Value = 3014.889
$ws.Range("I110").Value = 1428.2
$ws.Range("J110").Value = 4998.25
$ws.Range("K110").Value = 1428.2
$ws.Range("L110").Value = 4998.25
$ws.Range("M110").Value = 616.8
$ws.Range("N110").Value = -9088.25

$ws.Range("H122").Value = 2080.5833
$ws.Range("I122").Value = 1230.1111
$ws.Range("K122").Value = 3690.3333
$ws.Range("M122").Value = -1240.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 112.5
$ws.Range("I22").Value = 100
$ws.Range("J22").Value = 125
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = 125
$ws.Range("M22").Value = 73
$ws.Range("N22").Value = -471

$ws.Range("H82").Value = 22183.8
$ws.Range("I82").Value = 18275.7
$ws.Range("K82").Value = 18275.7
$ws.Range("M82").Value = -17892.7

$ws.Range("H85").Value = 22183.8
$ws.Range("I85").Value = 18275.7
$ws.Range("K85").Value = 18275.7
$ws.Range("M85").Value = -16949.7

$ws.Range("H105").Value = 6809.75
$ws.Range("I105").Value = 6413
$ws.Range("K105").Value = 6413
$ws.Range("M105").Value = -4666

$ws.Range("H132").Value = 230780
$ws.Range("J132").Value = 230780
$ws.Range("L132").Value = 230780
$ws.Range("N132").Value = -240900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 30000
$ws.Range("J53").Value = 30000
$ws.Range("L53").Value = 30000
$ws.Range("N53").Value = -31214

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3583
$ws.Range("I80").Value = 3374.5
$ws.Range("K80").Value = 3374.5
$ws.Range("M80").Value = -2376.5

$ws.Range("H83").Value = 3583
$ws.Range("I83").Value = 3374.5
$ws.Range("K83").Value = 16872.5
$ws.Range("M83").Value = -11880.5

$ws.Range("H123").Value = 28000.5
$ws.Range("J123").Value = 28000.5
$ws.Range("L123").Value = 28000.5
$ws.Range("N123").Value = -32900.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2325
$ws.Range("I16").Value = 1800
$ws.Range("K16").Value = 1800
$ws.Range("M16").Value = -1630

$ws.Range("H22").Value = 698.7143
$ws.Range("I22").Value = 539
$ws.Range("J22").Value = 1098
$ws.Range("K22").Value = 539
$ws.Range("L22").Value = 1098
$ws.Range("M22").Value = -244
$ws.Range("N22").Value = -1688

$ws.Range("H27").Value = 698.7143
$ws.Range("I27").Value = 539
$ws.Range("J27").Value = 1098
$ws.Range("K27").Value = 539
$ws.Range("L27").Value = 1098
$ws.Range("M27").Value = -432
$ws.Range("N27").Value = -1312

$ws.Range("H93").Value = 3187.5
$ws.Range("I93").Value = 3187.5
$ws.Range("K93").Value = 3187.5
$ws.Range("M93").Value = -1939.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 31363.545
$ws.Range("J26").Value = 31363.545
$ws.Range("L26").Value = 31363.545
$ws.Range("N26").Value = -31949.545

$ws.Range("H81").Value = 1267
$ws.Range("J81").Value = 1111
$ws.Range("L81").Value = 2222
$ws.Range("N81").Value = -4344

$ws.Range("H84").Value = 1267
$ws.Range("J84").Value = 1111
$ws.Range("L84").Value = 11110
$ws.Range("N84").Value = -21718

$ws.Range("H136").Value = 3840
$ws.Range("I136").Value = 3925
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 11775
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -9225
$ws.Range("N136").Value = -15600
